# Similar to "emph table": collapse the "Number of participants" /
# "Gender Male" / "Gender Female" rows into a single "Sex (male, %)" row,
# refresh the N= counts in the column headers, and update several
# statistics that were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-redundant rows ("Gender Male" / "Gender Female") -
# their data is folded into the (now renamed) row 2 "Sex (male, %)".
$ws.Rows("3:4").Delete()

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "Characteristics"
$ws.Range("B1").Value = "All participants (N=358)"
$ws.Range("C1").Value = "Participants with low BMI (N=178)"
$ws.Range("D1").Value = "Participants with high BMI (N=180)"
$ws.Range("E1").Value = "P value"

# --- Row 2: Sex (male, %) ----------------------------------------------
$ws.Range("A2").Value = "Sex (male, %)"
$ws.Range("B2").Value = "67 (19.0%)"
$ws.Range("C2").Value = "24 (13.5%)"
$ws.Range("D2").Value = "43 (24.0%)"
$ws.Range("E2").Value = 0.017

# --- Row 3: Age -----------------------------------------------------------
$ws.Range("A3").Value = "Age"
$ws.Range("B3").Value = "59.8±10.3"
$ws.Range("C3").Value = "61.3±10.8"
$ws.Range("D3").Value = "58.2±9.5"
$ws.Range("E3").Value = 0.003

# --- Row 4: Weight --------------------------------------------------------
$ws.Range("A4").Value = "Weight"
$ws.Range("B4").Value = "85.9±32.9"
$ws.Range("C4").Value = "55.0±6.3"
$ws.Range("D4").Value = "116.5±15.2"
$ws.Range("E4").Value = "<0.001"

# --- Row 5: Height --------------------------------------------------------
$ws.Range("A5").Value = "Height"
$ws.Range("B5").Value = "170.9±8.6"
$ws.Range("C5").Value = "171.2±8.3"
$ws.Range("D5").Value = "170.6±8.8"
$ws.Range("E5").Value = 0.534

# --- Row 6: Never smoker ---------------------------------------------------
$ws.Range("A6").Value = "Never smoker"
$ws.Range("B6").Value = "138 (39.0%)"
$ws.Range("C6").Value = "82 (46.1%)"
$ws.Range("D6").Value = "56 (30.0%)"
$ws.Range("E6").Value = 0.023

# --- Row 7: Ever smoker (no P value) ---------------------------------------
$ws.Range("A7").Value = "Ever smoker"
$ws.Range("B7").Value = "189 (53.0%)"
$ws.Range("C7").Value = "87 (48.9%)"
$ws.Range("D7").Value = "102 (57.0%)"
$ws.Range("E7").ClearContents()

# --- Row 8: Pack years ------------------------------------------------------
$ws.Range("A8").Value = "Pack years"
$ws.Range("B8").Value = "15.1±14.2"
$ws.Range("C8").Value = "15.2±13.2"
$ws.Range("D8").Value = "15.0±14.9"
$ws.Range("E8").Value = 0.947
